$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row number, new DAMSLTag (col I), new DialogAct (col J)
$changes = @(
    @(7, 'ba', 'Appreciation'),
    @(9, 'sd', 'Statement-non-opinion'),
    @(22, 'sd', 'Statement-non-opinion'),
    @(32, 'sd', 'Statement-non-opinion'),
    @(38, 'sd', 'Statement-non-opinion'),
    @(41, 'ba', 'Appreciation'),
    @(46, 'sd', 'Statement-non-opinion'),
    @(63, 'sd', 'Statement-non-opinion'),
    @(65, 'ba', 'Appreciation'),
    @(66, 'aa', 'Agree/Accept'),
    @(70, 'sd', 'Statement-non-opinion'),
    @(72, 'b', 'Acknowledge (Backchannel)'),
    @(77, 'b', 'Acknowledge (Backchannel)'),
    @(80, 'sv', 'Statement-opinion'),
    @(96, 'sv', 'Statement-opinion'),
    @(100, 'b', 'Acknowledge (Backchannel)'),
    @(102, 'b', 'Acknowledge (Backchannel)'),
    @(103, 'b', 'Acknowledge (Backchannel)'),
    @(107, 'sd', 'Statement-non-opinion'),
    @(108, '%', 'Uninterpretable'),
    @(110, 'sd', 'Statement-non-opinion'),
    @(114, 'b', 'Acknowledge (Backchannel)'),
    @(136, 'sd', 'Statement-non-opinion'),
    @(140, 'ba', 'Appreciation'),
    @(167, 'ba', 'Appreciation'),
    @(191, 'sv', 'Statement-opinion'),
    @(196, 'b', 'Acknowledge (Backchannel)'),
    @(202, 'sd', 'Statement-non-opinion'),
    @(205, 'sd', 'Statement-non-opinion'),
    @(206, 'aa', 'Agree/Accept'),
    @(208, 'aa', 'Agree/Accept'),
    @(209, 'b', 'Acknowledge (Backchannel)'),
    @(210, 'sd', 'Statement-non-opinion'),
    @(220, 'sv', 'Statement-opinion'),
    @(225, 'b', 'Acknowledge (Backchannel)'),
    @(228, 'b', 'Acknowledge (Backchannel)'),
    @(246, 'sd', 'Statement-non-opinion'),
    @(271, 'aa', 'Agree/Accept'),
    @(280, 'aa', 'Agree/Accept'),
    @(290, 'b', 'Acknowledge (Backchannel)'),
    @(298, 'b', 'Acknowledge (Backchannel)'),
    @(318, 'b', 'Acknowledge (Backchannel)'),
    @(330, 'ba', 'Appreciation'),
    @(333, 'ba', 'Appreciation'),
    @(335, 'ba', 'Appreciation'),
    @(337, 'ba', 'Appreciation'),
    @(340, 'ba', 'Appreciation'),
    @(344, 'b', 'Acknowledge (Backchannel)'),
    @(347, 'b', 'Acknowledge (Backchannel)'),
    @(363, 'aa', 'Agree/Accept'),
    @(368, 'b', 'Acknowledge (Backchannel)'),
    @(382, 'ba', 'Appreciation'),
    @(398, 'b', 'Acknowledge (Backchannel)'),
    @(402, 'b', 'Acknowledge (Backchannel)'),
    @(411, 'sv', 'Statement-opinion'),
    @(425, 'sv', 'Statement-opinion'),
    @(439, 'ba', 'Appreciation')
)

foreach ($change in $changes) {
    $row = $change[0]
    $damslTag = $change[1]
    $dialogAct = $change[2]
    $ws.Cells.Item($row, 9).Value = $damslTag
    $ws.Cells.Item($row, 10).Value = $dialogAct
}

Write-Output "Updated $($changes.Count) rows."